$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting rows 14:134 down to 15:135
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new record's data. Columns A, B, C, E,
# F, G, H, I, N, Q, R hold the same constant values as every other row in
# this dataset.
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(14, 3).Value = "La Araucanía"
$ws.Cells.Item(14, 4).Value = 44819
$ws.Cells.Item(14, 5).Value = 9
$ws.Cells.Item(14, 6).Value = 100112035
$ws.Cells.Item(14, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 55
$ws.Cells.Item(14, 11).Value = 24000
$ws.Cells.Item(14, 12).Value = 24000
$ws.Cells.Item(14, 13).Value = 24000
$ws.Cells.Item(14, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(14, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(14, 16).Value = 2400
$ws.Cells.Item(14, 17).Value = 10
$ws.Cells.Item(14, 18).Value = "Hortaliza"
